$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells AD1:AF1, copying the style/format of an existing
# header cell (AC1) so they match the rest of the header row (bold font,
# centered/top aligned, thin border), then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row (2-54)
# with the same values for each row.
$ws.Range("AD2:AD54").Value = 84
$ws.Range("AE2:AE54").Value = 78
$ws.Range("AF2:AF54").Value = 0
